$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ECS
$ws2 = $wb.Worksheets.Item(2)   # CSP

# --- Sheet1 (ECS): flip H3 and H5 Execute flags from N to Y ---
$ws1.Range("H3").Value = "Y"
$ws1.Range("H5").Value = "Y"

# --- Sheet1 (ECS): append new API test rows 7-9 ---
# Row 7
$ws1.Range("A7").Value = "API"
$ws1.Range("B7").Value = "Get Data API"
$ws1.Range("C7").Value = "tests.api.SampleAPI"
$ws1.Range("D7").Value = "'006"
$ws1.Range("E7").Value = "TC01_GetSomeData"
$ws1.Range("F7").Value = "env,browser"
$ws1.Range("G7").Value = "ECS_API,API"
$ws1.Range("H7").Value = "Y"
$ws1.Range("I7").Value = "null"

# Row 8
$ws1.Range("A8").Value = "API"
$ws1.Range("B8").Value = "ReqRes API"
$ws1.Range("C8").Value = "tests.api.SampleAPI"
$ws1.Range("D8").Value = "'007"
$ws1.Range("E8").Value = "TC02_PostSomeData"
$ws1.Range("F8").Value = "env,browser"
$ws1.Range("G8").Value = "ECS_API,API"
$ws1.Range("H8").Value = "Y"
$ws1.Range("I8").Value = "null"

# Row 9
$ws1.Range("A9").Value = "API"
$ws1.Range("B9").Value = "ReqRes API"
$ws1.Range("C9").Value = "tests.api.SampleAPI"
$ws1.Range("D9").Value = "'008"
$ws1.Range("E9").Value = "TC03_PutSomeData"
$ws1.Range("F9").Value = "env,browser"
$ws1.Range("G9").Value = "ECS_API,API"
$ws1.Range("H9").Value = "Y"
$ws1.Range("I9").Value = "TC02_PostSomeData"

# --- update the view/selection on ECS, then restore CSP as the active sheet ---
$null = $ws1.Range("I7").Select()
$null = $ws2.Activate()
$null = $ws2.Range("H3").Select()
